$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36; this shifts rows 36..79 down to 37..80
# and extends the used range to row 80 (matching <dimension ref="A1:R80"/>).
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with the new weekly price record.
$ws.Range("A36").Value = 10
$ws.Range("B36").Value = "Vega Modelo de Temuco"
$ws.Range("C36").Value = "La Araucanía"
$ws.Range("D36").Value = 44483
$ws.Range("E36").Value = 9
$ws.Range("F36").Value = 100112031
$ws.Range("G36").Value = "Poroto verde"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 200
$ws.Range("K36").Value = 3000
$ws.Range("L36").Value = 3000
$ws.Range("M36").Value = 3000
$ws.Range("N36").Value = "$/kilo"
$ws.Range("O36").Value = "Provincia de Limarí"
$ws.Range("P36").Value = 3000
$ws.Range("Q36").Value = 1
$ws.Range("R36").Value = "Hortaliza"
